$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New shared strings must be created in this exact order so the
# resulting sharedStrings.xml table matches the authored workbook.
$ws.Range("B60").Value = "Trains off, worked from uni on writing disseration - primarily doing captioning"
$ws.Range("B58").Value = "Writing dissertation - first half of Network analyiser section done"
$ws.Range("B59").Value = "Writing dissertation - second half of Network analyiser and most of topic modeller done"
$ws.Range("B56").Value = "Supervision with Gabrila and then worked from the uni"
$ws.Range("B57").Value = "AL"

# B56 sits in a "week boundary" row (same formatting family as B6, B11,
# B16, ... B51) which carries a thin bottom border.
$ws.Range("B56").Borders.Item(9).LineStyle = 1

$ws.Range("B57").Select()
